# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (want-to-go count) figures for the two anime-expo
# entries that appear on both the "展览" (exhibitions) sheet and the
# "全部类型" (all types) roll-up sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1298   # 南宁·熊喵M动漫嘉年华【免费】   1295 -> 1298
$ws1.Range("F3").Value = 2810   # 南宁·第二届北极光动漫展        2806 -> 2810
$ws1.Range("F4").Value = 253    # 南宁·万圣漫控嘉年华10          252  -> 253

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1298   # 南宁·熊喵M动漫嘉年华【免费】   1295 -> 1298
$ws4.Range("F4").Value = 2810   # 南宁·第二届北极光动漫展        2806 -> 2810
$ws4.Range("F6").Value = 253    # 南宁·万圣漫控嘉年华10          252  -> 253
